# Regenerate Report for Handback
#
# The handback-status report is re-generated by the CI job; one of the two
# handed-back source files (f2bfd839-ba30-4254-a8da-68a9fd3cc98c.*) picked up
# fresher "xliff generate" / handoff / handback timestamps on this run, while
# the other file (ac048824-46ca-4c8f-a70c-1b15639240f9.*) is untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
# Row 3 is the f2bfd839... file; its "Latest HO Xliff Generate Date" column
# (G) moves forward to the new generation timestamp. Row 2 (ac048824...) is
# left as-is.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-19 12:47:11"

# --- zh-cn sheet -------------------------------------------------------------
# Row 3 (f2bfd839...): Correspond Handoff Datetime (H) and Correspond
# Handback DateTime (K) both get refreshed timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-19 12:47:02"
$wsZhCn.Range("K3").Value = "2016-08-19 12:47:32"

# --- de-de sheet -------------------------------------------------------------
# Same story for the de-de table: row 3 (f2bfd839...) gets new Handoff /
# Handback datetimes.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-19 12:47:11"
$wsDeDe.Range("K3").Value = "2016-08-19 12:47:39"
